# Added data for ProvarCache
# Updates the "RMA Details Maintenance Grid" sheet's RMA/line-item references
# to the freshly generated RMA-DU0X batch (replacing the prior RMA-EM0X batch).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

$ws.Range("E2").Value = "RMA-DU0X-001"
$ws.Range("F2").Value = "RMA-DU0X-1-1"
$ws.Range("J2").Value = "a7s5f000000xKMAAA2"

$ws.Range("E3").Value = "RMA-DU0X-002"
$ws.Range("F3").Value = "RMA-DU0X-1-2"
$ws.Range("J3").Value = "a7s5f000000xKMBAA2"

$ws.Range("E4").Value = "RMA-DU0X-003"
$ws.Range("F4").Value = "RMA-DU0X-1-3"
$ws.Range("J4").Value = "a7s5f000000xKMCAA2"
